$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7463439233484619
$ws.Range("C2").Value = 0.6281833616298812
$ws.Range("D2").Value = 0.6821848352154875

$ws.Range("B3").Value = 0.4352030947775629
$ws.Range("C3").Value = 0.5730050933786078
$ws.Range("D3").Value = 0.4946866984243313

$ws.Range("B4").Value = 0.6097906055461234
$ws.Range("C4").Value = 0.6097906055461234
$ws.Range("D4").Value = 0.6097906055461234
$ws.Range("E4").Value = 0.6097906055461234

$ws.Range("B5").Value = 0.5907735090630124
$ws.Range("C5").Value = 0.6005942275042445
$ws.Range("D5").Value = 0.5884357668199094

$ws.Range("B6").Value = 0.6426303138248289
$ws.Range("C6").Value = 0.6097906055461234
$ws.Range("D6").Value = 0.6196854562851021
